$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'256.12"
$ws.Range("E2").Value = "'0.19%"
$ws.Range("G2").Value = "'14"
$ws.Range("D3").Value = "'27.10"
$ws.Range("E3").Value = "'-4.44%"
$ws.Range("G3").Value = "'14"
$ws.Range("D4").Value = "'4.625"
$ws.Range("E4").Value = "'-10.82%"
$ws.Range("G4").Value = "'14"
$ws.Range("D5").Value = "'0.05891"
$ws.Range("E5").Value = "'0.57%"
$ws.Range("G5").Value = "'14"
$ws.Range("D6").Value = "'6.640"
$ws.Range("E6").Value = "'-1.24%"
$ws.Range("G6").Value = "'14"
$ws.Range("D7").Value = "'0.8686"
$ws.Range("E7").Value = "'-0.16%"
$ws.Range("G7").Value = "'14"
$ws.Range("D8").Value = "'0.9448"
$ws.Range("E8").Value = "'-1.52%"
$ws.Range("G8").Value = "'14"
$ws.Range("B9").Value = "WazirX"
$ws.Range("C9").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D9").Value = "'0.1404"
$ws.Range("E9").Value = "'-0.45%"
$ws.Range("G9").Value = "'14"
$ws.Range("B10").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C10").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D10").Value = "'0.03728"
$ws.Range("E10").Value = "'8.14%"
$ws.Range("G10").Value = "'14"
$ws.Range("B11").Value = "MandalaExchangeToken"
$ws.Range("C11").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D11").Value = "'0.07101"
$ws.Range("E11").Value = "'-0.71%"
$ws.Range("G11").Value = "'14"
$ws.Range("B12").Value = "BitrueCoin"
$ws.Range("C12").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D12").Value = "'0.03200"
$ws.Range("E12").Value = "'-0.25%"
$ws.Range("G12").Value = "'14"
$ws.Range("B13").Value = "BitMartToken"
$ws.Range("C13").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D13").Value = "'0.09262"
$ws.Range("E13").Value = "'0.43%"
$ws.Range("G13").Value = "'14"
$ws.Range("B14").Value = "BitForexToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D14").Value = "'0.001538"
$ws.Range("E14").Value = "'-0.43%"
$ws.Range("G14").Value = "'14"
$ws.Range("B15").Value = "One"
$ws.Range("C15").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D15").Value = "'0.0006004"
$ws.Range("E15").Value = "'-1.06%"
$ws.Range("G15").Value = "'14"
$ws.Range("D16").Value = "'0.006133"
$ws.Range("E16").Value = "'2.16%"
$ws.Range("G16").Value = "'14"
$ws.Range("D17").Value = "'3.512"
$ws.Range("E17").Value = "'0.37%"
$ws.Range("G17").Value = "'14"
$ws.Range("E18").Value = "'-1.27%"
$ws.Range("G18").Value = "'14"
$ws.Range("D19").Value = "'2.218"
$ws.Range("E19").Value = "'-0.32%"
$ws.Range("G19").Value = "'14"
$ws.Range("D20").Value = "'0.3111"
$ws.Range("E20").Value = "'-2.10%"
$ws.Range("G20").Value = "'14"
$ws.Range("E21").Value = "'-1.20%"
$ws.Range("G21").Value = "'14"
$ws.Range("D22").Value = "'3.849"
$ws.Range("E22").Value = "'9.11%"
$ws.Range("G22").Value = "'14"
$ws.Range("D23").Value = "'0.04235"
$ws.Range("E23").Value = "'1.22%"
$ws.Range("G23").Value = "'14"
$ws.Range("E24").Value = "'-1.32%"
$ws.Range("G24").Value = "'14"
$ws.Range("D25").Value = "'0.001223"
$ws.Range("E25").Value = "'0.20%"
$ws.Range("G25").Value = "'14"
$ws.Range("D26").Value = "'0.004281"
$ws.Range("E26").Value = "'-6.27%"
$ws.Range("G26").Value = "'14"
$ws.Range("D27").Value = "'0.0001200"
$ws.Range("E27").Value = "'0.01%"
$ws.Range("G27").Value = "'14"
$ws.Range("D28").Value = "'0.0001503"
$ws.Range("E28").Value = "'2.50%"
$ws.Range("G28").Value = "'14"
$ws.Range("G29").Value = "'14"
$ws.Range("G30").Value = "'14"
$ws.Range("G31").Value = "'14"
$ws.Range("G32").Value = "'14"
$ws.Range("G33").Value = "'14"
$ws.Range("G34").Value = "'14"
$ws.Range("G35").Value = "'14"
$ws.Range("G36").Value = "'14"
$ws.Range("G37").Value = "'14"
$ws.Range("G38").Value = "'14"
$ws.Range("G39").Value = "'14"
$ws.Range("D40").Value = "'0.03812"
$ws.Range("E40").Value = "'-0.24%"
$ws.Range("G40").Value = "'14"
$ws.Range("D41").Value = "'0.006240"
$ws.Range("E41").Value = "'11.04%"
$ws.Range("G41").Value = "'14"
$ws.Range("E42").Value = "'-0.32%"
$ws.Range("G42").Value = "'14"
$ws.Range("D43").Value = "'0.002430"
$ws.Range("E43").Value = "'3.69%"
$ws.Range("G43").Value = "'14"
$ws.Range("D44").Value = "'0.01121"
$ws.Range("E44").Value = "'15.16%"
$ws.Range("G44").Value = "'14"
$ws.Range("D45").Value = "'0.00005504"
$ws.Range("E45").Value = "'1.78%"
$ws.Range("G45").Value = "'14"
$ws.Range("D46").Value = "'0.00000000751"
$ws.Range("E46").Value = "'0.08%"
$ws.Range("G46").Value = "'14"
$ws.Range("D47").Value = "'0.07786"
$ws.Range("E47").Value = "'-13.47%"
$ws.Range("G47").Value = "'14"
$ws.Range("D48").Value = "'0.002280"
$ws.Range("E48").Value = "'7.20%"
$ws.Range("G48").Value = "'14"
$ws.Range("D49").Value = "'0.00002102"
$ws.Range("E49").Value = "'0.08%"
$ws.Range("G49").Value = "'14"
$ws.Range("E50").Value = "'0.08%"
$ws.Range("G50").Value = "'14"
$ws.Range("G51").Value = "'14"
